# Update "Fecha" (D), "Volumen" (M), "Precio mínimo" (N), "Precio máximo" (O),
# "Precio promedio ponderado" (P) and "Precio $/Kg" (S) values for rows 2-9
# of the active worksheet to match the new weekly data snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    3 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    4 = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    5 = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
    6 = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    7 = @{ D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 }
    8 = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    9 = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # Column D - Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M   # Column M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # Column N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # Column O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # Column P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals.S   # Column S - Precio $/Kg
}
